$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 305.44446
$ws.Range("I6").Value = 268.625
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 805.875
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -693.875
$ws.Range("N6").Value = -2024
$ws.Range("H33").Value = 246.89189
$ws.Range("I33").Value = 201.02942
$ws.Range("K33").Value = 201.02942
$ws.Range("M33").Value = 27.97058000000001
$ws.Range("H64").Value = 33101.94
$ws.Range("I64").Value = 3663.3
$ws.Range("J64").Value = 75157.14
$ws.Range("K64").Value = 3663.3
$ws.Range("L64").Value = 75157.14
$ws.Range("M64").Value = -3415.3
$ws.Range("N64").Value = -75653.14
$ws.Range("H67").Value = 33101.94
$ws.Range("I67").Value = 3663.3
$ws.Range("J67").Value = 75157.14
$ws.Range("K67").Value = 3663.3
$ws.Range("L67").Value = 75157.14
$ws.Range("M67").Value = -2805.3
$ws.Range("N67").Value = -76873.14
$ws.Range("H86").Value = 33336404
$ws.Range("I86").Value = 2085.9048
$ws.Range("J86").Value = 111116480
$ws.Range("K86").Value = 2085.9048
$ws.Range("L86").Value = 111116480
$ws.Range("M86").Value = -962.9047999999998
$ws.Range("N86").Value = -111118726
$ws.Range("H89").Value = 33336404
$ws.Range("I89").Value = 2085.9048
$ws.Range("J89").Value = 111116480
$ws.Range("K89").Value = 10429.524
$ws.Range("L89").Value = 555582400
$ws.Range("M89").Value = -4813.523999999999
$ws.Range("N89").Value = -555593632
$ws.Range("H120").Value = 40761
$ws.Range("J120").Value = 40761
$ws.Range("L120").Value = 40761
$ws.Range("N120").Value = -50437
$ws.Range("H137").Value = 23810786
$ws.Range("I137").Value = 1171.9722
$ws.Range("J137").Value = 166668460
$ws.Range("K137").Value = 3515.9166
$ws.Range("L137").Value = 500005380
$ws.Range("M137").Value = -965.9165999999996
$ws.Range("N137").Value = -500010480
$ws.Range("H141").Value = 720.25
$ws.Range("I141").Value = 720.25
$ws.Range("K141").Value = 2160.75
$ws.Range("M141").Value = 3019.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2641.25
$ws.Range("I61").Value = 1947.1428
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 1947.1428
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -1735.1428
$ws.Range("N61").Value = -7924
$ws.Range("H74").Value = 4834.0713
$ws.Range("I74").Value = 866.5263
$ws.Range("J74").Value = 13210
$ws.Range("K74").Value = 866.5263
$ws.Range("L74").Value = 13210
$ws.Range("M74").Value = 7.473700000000008
$ws.Range("N74").Value = -14958
$ws.Range("H77").Value = 4834.0713
$ws.Range("I77").Value = 866.5263
$ws.Range("J77").Value = 13210
$ws.Range("K77").Value = 4332.6315
$ws.Range("L77").Value = 66050
$ws.Range("M77").Value = 35.36850000000049
$ws.Range("N77").Value = -74786
$ws.Range("H132").Value = 1718.5758
$ws.Range("I132").Value = 1692.0741
$ws.Range("J132").Value = 1837.8334
$ws.Range("K132").Value = 5076.2223
$ws.Range("L132").Value = 5513.5002
$ws.Range("M132").Value = -2546.2223
$ws.Range("N132").Value = -10573.5002
$ws.Range("H136").Value = 2641.25
$ws.Range("I136").Value = 1947.1428
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 5841.428400000001
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -3291.428400000001
$ws.Range("N136").Value = -27600
$ws.Range("H139").Value = 59933.332
$ws.Range("J139").Value = 59933.332
$ws.Range("L139").Value = 59933.332
$ws.Range("N139").Value = -70213.33199999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("N132").Value = 0
$ws.Range("H134").Value = 59898
$ws.Range("I134").Value = 63585.375
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 190756.125
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -188221.125
$ws.Range("N134").Value = -7770
$ws.Range("H137").Value = 95770
$ws.Range("J137").Value = 95770
$ws.Range("L137").Value = 95770
$ws.Range("N137").Value = -105970
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2472
$ws.Range("I58").Value = 1681.4
$ws.Range("J58").Value = 3190.7273
$ws.Range("K58").Value = 1681.4
$ws.Range("L58").Value = 3190.7273
$ws.Range("M58").Value = -1478.4
$ws.Range("N58").Value = -3596.7273
$ws.Range("H62").Value = 2999.5
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4247
$ws.Range("H65").Value = 2999.5
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21235
$ws.Range("H132").Value = 4916.933
$ws.Range("I132").Value = 4161.4546
$ws.Range("J132").Value = 6994.5
$ws.Range("K132").Value = 12484.3638
$ws.Range("L132").Value = 20983.5
$ws.Range("M132").Value = -9954.363799999999
$ws.Range("N132").Value = -26043.5
$ws.Range("H134").Value = 3361.9524
$ws.Range("I134").Value = 3881.375
$ws.Range("J134").Value = 1699.8
$ws.Range("K134").Value = 11644.125
$ws.Range("L134").Value = 5099.4
$ws.Range("M134").Value = -9109.125
$ws.Range("N134").Value = -10169.4
$ws.Range("H136").Value = 2472
$ws.Range("I136").Value = 1681.4
$ws.Range("J136").Value = 3190.7273
$ws.Range("K136").Value = 5044.200000000001
$ws.Range("L136").Value = 9572.1819
$ws.Range("M136").Value = -2494.200000000001
$ws.Range("N136").Value = -14672.1819

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 525
$ws.Range("I99").Value = 525
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1575
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 671
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("H129").Value = 936.5
$ws.Range("J129").Value = 1118.2307
$ws.Range("L129").Value = 3354.6921
$ws.Range("N129").Value = -13354.6921
$ws.Range("H132").Value = 76924860
$ws.Range("I132").Value = 200001100
$ws.Range("J132").Value = 2212.5
$ws.Range("K132").Value = 1800009900
$ws.Range("L132").Value = 19912.5
$ws.Range("M132").Value = -1800007370
$ws.Range("N132").Value = -24972.5
$ws.Range("H133").Value = 8753.75
$ws.Range("I133").Value = 5015
$ws.Range("K133").Value = 15045
$ws.Range("M133").Value = -9985

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 82.86957
$ws.Range("I2").Value = 25.285715
$ws.Range("K2").Value = 25.285715
$ws.Range("M2").Value = 87.714285
$ws.Range("H132").Value = 2166.4285
$ws.Range("I132").Value = 1630.4
$ws.Range("J132").Value = 3506.5
$ws.Range("K132").Value = 4891.200000000001
$ws.Range("L132").Value = 10519.5
$ws.Range("M132").Value = -2361.200000000001
$ws.Range("N132").Value = -15579.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1038.409
$ws.Range("I136").Value = 902.9706
$ws.Range("J136").Value = 1498.9
$ws.Range("K136").Value = 2708.9118
$ws.Range("L136").Value = 4496.700000000001
$ws.Range("M136").Value = -158.9117999999999
$ws.Range("N136").Value = -9596.700000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2244.0444
$ws.Range("I132").Value = 2150.7878
$ws.Range("J132").Value = 2500.5
$ws.Range("K132").Value = 6452.3634
$ws.Range("L132").Value = 7501.5
$ws.Range("M132").Value = -3922.3634
$ws.Range("N132").Value = -12561.5
